$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row in column A (data starts at row 2, header at row 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $typeCell = $ws.Cells.Item($r, 2)
    $valueCell = $ws.Cells.Item($r, 3)

    $typeVal = $typeCell.Value2
    $valueVal = $valueCell.Value2

    if ([string]::IsNullOrEmpty($typeVal)) { continue }

    # Determine gender/coed suffix from the existing value column text
    $suffix = $null
    if ($valueVal -like "*Boys*") { $suffix = "boys" }
    elseif ($valueVal -like "*Girls*") { $suffix = "girls" }
    elseif ($valueVal -like "*Coed*") { $suffix = "coed" }

    if ($typeVal -eq "club-sports") {
        $typeCell.Value = "sports_club_$suffix"
    } elseif ($typeVal -eq "uil-sports") {
        $typeCell.Value = "sports_uil_$suffix"
    }

    $valueCell.Value = "Track/Field"
}
